$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45918.01041666666, 1683),
    @(45918.02083333334, 1693),
    @(45918.03125, 1697),
    @(45918.04166666666, 1703),
    @(45918.05208333334, 1845),
    @(45918.0625, 1871),
    @(45918.07291666666, 1883),
    @(45918.08333333334, 1898),
    @(45918.09375, 2004),
    @(45918.10416666666, 2007),
    @(45918.11458333334, 2007),
    @(45918.125, 2031),
    @(45918.13541666666, 2048),
    @(45918.14583333334, 2051),
    @(45918.15625, 2045),
    @(45918.16666666666, 2041),
    @(45918.17708333334, 2059),
    @(45918.1875, 2059),
    @(45918.19791666666, 2061),
    @(45918.20833333334, 2065),
    @(45918.21875, 2074),
    @(45918.22916666666, 2077),
    @(45918.23958333334, 2068),
    @(45918.25, 2070),
    @(45918.26041666666, 2084),
    @(45918.27083333334, 2089),
    @(45918.28125, 2091),
    @(45918.29166666666, 2080),
    @(45918.30208333334, 2030),
    @(45918.3125, 2034),
    @(45918.32291666666, 2035),
    @(45918.33333333334, 2038),
    @(45918.34375, 2213),
    @(45918.35416666666, 2223),
    @(45918.36458333334, 2230),
    @(45918.375, 2240),
    @(45918.38541666666, 2351),
    @(45918.39583333334, 2355),
    @(45918.40625, 2358),
    @(45918.41666666666, 2361),
    @(45918.42708333334, 2346),
    @(45918.4375, 2347),
    @(45918.44791666666, 2347),
    @(45918.45833333334, 2348),
    @(45918.46875, 2363),
    @(45918.47916666666, 2362),
    @(45918.48958333334, 2362),
    @(45918.5, 2361),
    @(45918.51041666666, 2352),
    @(45918.52083333334, 2351),
    @(45918.53125, 2351),
    @(45918.54166666666, 2350),
    @(45918.55208333334, 2220),
    @(45918.5625, 2219),
    @(45918.57291666666, 2219),
    @(45918.58333333334, 2219),
    @(45918.59375, 2149),
    @(45918.60416666666, 2147),
    @(45918.61458333334, 2145),
    @(45918.625, 2143),
    @(45918.63541666666, 2237),
    @(45918.64583333334, 2233),
    @(45918.65625, 2229),
    @(45918.66666666666, 2225),
    @(45918.67708333334, 1971),
    @(45918.6875, 1962),
    @(45918.69791666666, 1952),
    @(45918.70833333334, 1942),
    @(45918.71875, 1556),
    @(45918.72916666666, 1540),
    @(45918.73958333334, 1525),
    @(45918.75, 1509),
    @(45918.76041666666, 1246),
    @(45918.77083333334, 1240),
    @(45918.78125, 1234),
    @(45918.79166666666, 1228),
    @(45918.80208333334, 1068),
    @(45918.8125, 1065),
    @(45918.82291666666, 1063),
    @(45918.83333333334, 1061),
    @(45918.84375, 995),
    @(45918.85416666666, 994),
    @(45918.86458333334, 993),
    @(45918.875, 993),
    @(45918.88541666666, 987),
    @(45918.89583333334, 987),
    @(45918.90625, 987),
    @(45918.91666666666, 988),
    @(45918.92708333334, 1038),
    @(45918.9375, 1039),
    @(45918.94791666666, 1041),
    @(45918.95833333334, 1044),
    @(45918.96875, 0),
    @(45918.97916666666, 0),
    @(45918.98958333334, 0),
    @(45919.0, 0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
